$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (e.g. AC1) so the
# new header cells match the bold/centered/bordered look of the rest of
# row 1.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Team record values for every data row (2-40)
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 93
    $ws.Cells.Item($r, 31).Value = 69
    $ws.Cells.Item($r, 32).Value = 0
}
